# Insert a new daily price record as row 178 in the "Hortaliza, Femacal de
# La Calera - Poroto granado" sheet. Inserting a whole row shifts the
# existing rows 178-270 down to 179-271 (preserving their data/formatting)
# and extends the sheet dimension from R270 to R271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(178).Insert()

$ws.Cells.Item(178, 1).Value2  = 3
$ws.Cells.Item(178, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(178, 3).Value2  = "Coquimbo"
$ws.Cells.Item(178, 4).Value2  = 45001
$ws.Cells.Item(178, 5).Value2  = 5
$ws.Cells.Item(178, 6).Value2  = 100112030
$ws.Cells.Item(178, 7).Value2  = "Poroto granado"
$ws.Cells.Item(178, 8).Value2  = "Sin especificar"
$ws.Cells.Item(178, 9).Value2  = "Primera"
$ws.Cells.Item(178, 10).Value2 = 82
$ws.Cells.Item(178, 11).Value2 = 34000
$ws.Cells.Item(178, 12).Value2 = 35000
$ws.Cells.Item(178, 13).Value2 = 34488
$ws.Cells.Item(178, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(178, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(178, 16).Value2 = 1380
$ws.Cells.Item(178, 17).Value2 = 25
$ws.Cells.Item(178, 18).Value2 = "Hortaliza"
